# Finalized project and added VST3 and Standalone releases
# Updates the "Sprint 1 - Bilan" worksheet: marks most remaining tasks as
# "Finis" (done), fills in their progress (Avancement) and time invested
# (temp investi), and records two new delay/comment notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1 - Bilan")

# Row 7 (task 4 - "Créer un grain"): time invested Moyen -> Long
$ws.Range("F7").Value = "Long"

# Row 13 (task 10): was "En cours" -> now "Finis", 100% avancement, Moyen
$ws.Range("D13").Value = "Finis"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = "Moyen"

# Row 14 (task 11): was "Éventuel" -> now "Finis", 100% avancement, Moyen,
# plus a note about std::rand()
$ws.Range("D14").Value = "Finis"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = "Moyen"

# Row 15 (task 12): was "Éventuel" -> now "Finis", 100% avancement, Court
$ws.Range("D15").Value = "Finis"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = "Court"

# Row 16 (task 13): was "Éventuel" -> now "Finis", 100% avancement, Court
$ws.Range("D16").Value = "Finis"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = "Court"

# Row 17 (task 14): stays "Éventuel", but now explicitly 0% avancement, Null
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = "Null"

# Row 20 (task 17): time invested Long -> Court
$ws.Range("F20").Value = "Court"

# Row 21 (task 18): was "Éventuel" -> now "Finis", 100% avancement, Moyen
$ws.Range("D21").Value = "Finis"
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = "Moyen"

# Row 22 (task 19): was "Éventuel" -> now "Finis", 100% avancement, Très long,
# plus a note about VST parameter creation complications.
# (Set before row 14's comment so new shared strings are appended in the
# same order as the target workbook: "Complications..." then "std::rand()...")
$ws.Range("D22").Value = "Finis"
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = "Très long"
$ws.Range("G22").Value = "Complications lors de création de paramètre VST après construction"

$ws.Range("G14").Value = "std::rand() retourne une valeur maximum trop petite"

# Rows 23-29 (tasks 20,21,22,23,24,25,26): all were "Éventuel" -> now "Finis",
# 100% avancement, Court
$ws.Range("D23").Value = "Finis"
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = "Court"

$ws.Range("D24").Value = "Finis"
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = "Court"

$ws.Range("D25").Value = "Finis"
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = "Court"

$ws.Range("D26").Value = "Finis"
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = "Court"

$ws.Range("D27").Value = "Finis"
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = "Court"

$ws.Range("D28").Value = "Finis"
$ws.Range("E28").Value = 1
$ws.Range("F28").Value = "Court"

$ws.Range("D29").Value = "Finis"
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = "Court"

# Update the active selection to reflect where the author ended up working
$ws.Range("G29").Select()
